$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the training-part split - reuse the same header
# style (bold, centered, bordered) as the existing headers.
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Register the custom date/time number format (engine first records the
# lower-case attempt as numFmtId 164, then the final upper-case format as
# numFmtId 165 - apply on a single cell first so both format ids are
# registered but only the final one is left in use).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Timestamps (as Excel date serials) replacing the old inline-string
# timestamps, plus the new "Trening" column values for existing + newly
# added rows.
$data = @(
    @(2,  45685.64974571759, 1253,   14.87,              3.982890571866718, "10-15", "Duża Gra"),
    @(3,  45685.65119016204, 1377.8, 14.88,              3.912773438862392, "10-15", "Duża Gra"),
    @(4,  45685.66514965278, 2583.9, 14.19,              3.910629987716676, "10-15", "Duża Gra"),
    @(5,  45685.64706979167, 1021.8, 9.46,               3.141187412398201, "5-10",  "Duża Gra"),
    @(6,  45685.64974108796, 1252.6, 8.44,               2.973512342997959, "5-10",  "Duża Gra"),
    @(7,  45685.66514618055, 2583.6, 8.880000000000001,  3.02509641647339,  "5-10",  "Duża Gra"),
    @(8,  45685.671309375,   3116.1, 12.84,              3.168988500322613, "10-15", "Mała Gra"),
    @(9,  45685.6777630787,  3673.7, 13.99,              3.344213451657978, "10-15", "Mała Gra"),
    @(10, 45685.68089155092, 3944,   12.23,              3.53179751123701,  "10-15", "Mała Gra"),
    @(11, 45685.66746215278, 2783.7, 9.06,               3.013893025262014, "5-10",  "Mała Gra"),
    @(12, 45685.67130706018, 3115.9, 9.960000000000001,  2.853109666279381, "5-10",  "Mała Gra"),
    @(13, 45685.68088923611, 3943.8, 9.57,               3.297971963882447, "5-10",  "Mała Gra")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
